# Apply the centrality.xlsx update described in the commit:
# 1) Two pairs of tied-rank characters swap rows' worth of data (names + stats),
#    reflecting a reorder that happened upstream when the notebook was re-run.
# 2) A handful of floating point centrality values shift by one ULP due to a
#    (nondeterministic) recomputation of eigenvector_centrality / other stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "name" (column G) values for the two tied pairs of rows ---
$nameRow77  = $ws.Range("G77").Value()
$nameRow120 = $ws.Range("G120").Value()
$ws.Range("G77").Value  = $nameRow120
$ws.Range("G120").Value = $nameRow77

$nameRow85  = $ws.Range("G85").Value()
$nameRow125 = $ws.Range("G125").Value()
$ws.Range("G85").Value  = $nameRow125
$ws.Range("G125").Value = $nameRow85

# --- Swap the degree/degree_centrality/eigenvector_centrality/in_degree stats
#     for rows 77 <-> 120 (columns C, D, E, F) ---
$ws.Range("C77").Value = [double]"0"
$ws.Range("D77").Value = [double]"0"
$ws.Range("E77").Value = [double]"1.284276084573828e-13"
$ws.Range("F77").Value = [double]"0"

$ws.Range("C120").Value = [double]"1"
$ws.Range("D120").Value = [double]"0.008130081300813009"
$ws.Range("E120").Value = [double]"0.007242707713995649"
$ws.Range("F120").Value = [double]"1"

# --- Swap the degree/degree_centrality/out_degree stats for rows 85 <-> 125
#     (columns C, D, H) ---
$ws.Range("C85").Value = [double]"3"
$ws.Range("D85").Value = [double]"0.02439024390243903"
$ws.Range("H85").Value = [double]"3"

$ws.Range("C125").Value = [double]"5"
$ws.Range("D125").Value = [double]"0.04065040650406505"
$ws.Range("H125").Value = [double]"5"

# --- Tiny (last-digit) floating point adjustments to eigenvector_centrality
#     and a couple of betweenness_centrality values, matching the refreshed
#     notebook run ---
$ws.Range("E3").Value   = [double]"0.004318535650622628"
$ws.Range("E6").Value   = [double]"0.002984730610197574"
$ws.Range("E10").Value  = [double]"0.08066773161898984"
$ws.Range("E18").Value  = [double]"0.10827424760359"
$ws.Range("E23").Value  = [double]"0.3217690234418081"
$ws.Range("E27").Value  = [double]"0.05754444355901207"
$ws.Range("E36").Value  = [double]"0.2841119899212121"
$ws.Range("E39").Value  = [double]"0.03918421343610794"
$ws.Range("B42").Value  = [double]"0.004800646584853955"
$ws.Range("E44").Value  = [double]"0.0004005461788657446"
$ws.Range("B49").Value  = [double]"0.00559176014584846"
$ws.Range("E51").Value  = [double]"0.06673976267976568"
$ws.Range("E56").Value  = [double]"0.01368845259232841"
$ws.Range("E61").Value  = [double]"0.000824331796206526"
$ws.Range("E66").Value  = [double]"0.09477411901755346"
$ws.Range("E70").Value  = [double]"0.02196410317560709"
$ws.Range("E96").Value  = [double]"0.04336702136584433"
$ws.Range("E102").Value = [double]"0.1094806709571929"
$ws.Range("E105").Value = [double]"0.07247327824352919"
$ws.Range("E110").Value = [double]"0.2033701472773884"
$ws.Range("E117").Value = [double]"0.01326347069326225"
$ws.Range("E121").Value = [double]"0.02196410317560709"
$ws.Range("E122").Value = [double]"0.08354577196206246"
$ws.Range("B124").Value = [double]"0.03299247227728684"
$ws.Range("E124").Value = [double]"0.2824910616271694"
